# Wireshark Helper document edits
$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "has to" -> "must"  (== operator explanation)
# ------------------------------------------------------------------
$d.Content.Find.Execute("has to", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "must", 2) | Out-Null

# ------------------------------------------------------------------
# Helper: locate paragraphs by their (unique) current text so later
# edits are not dependent on brittle, hardcoded paragraph indices.
# ------------------------------------------------------------------
function Find-ParagraphByText($doc, $needle) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text -like $needle) {
            return $p
        }
    }
    return $null
}

# ------------------------------------------------------------------
# 2) "Example: ip.addr == 192.168.0.0" paragraph gets an explanation
#    appended (tab separated) describing what the filter returns.
# ------------------------------------------------------------------
$p1 = Find-ParagraphByText $d "Example: ip.addr == 192.168.0.0*"
$newText1 = "Example: ip.addr == 192.168.0.0 `t- Gives all packets where ip address 192.168.0.0 `t`t`t`t`t`tis source or destination"
$p1.Range.Text = $newText1

# ------------------------------------------------------------------
# 3) "- Gives " paragraph gets its explanation filled in.
# ------------------------------------------------------------------
$p2 = Find-ParagraphByText $d "- Gives `r"
if (-not $p2) { $p2 = Find-ParagraphByText $d "- Gives" }
$p2FirstLine = $p2.FirstLineIndent
$newText2 = "- Gives all packets where destination ip is `t`t`t`t`t`t`t192.168.0.0 and source ip is 63.179.30.136"
$p2.Range.Text = $newText2
$p2.Style = "Normal"
$p2.FirstLineIndent = $p2FirstLine

# ------------------------------------------------------------------
# 4) The bulleted "Gives the packets that have the ip address..."
#    paragraph loses its numbered-list formatting and instead gets a
#    literal "- " prefix with a fixed left indent.
# ------------------------------------------------------------------
$p3 = Find-ParagraphByText $d "Gives the packets that have the ip address*"
$p3.Style = "Normal"
$p3.LeftIndent = 216
$insEnd = $p3.Range.Start
$insRange = $d.Range($insEnd, $insEnd)
$insRange.InsertBefore("- ")
$newRunRange = $d.Range($p3.Range.Start, $p3.Range.Start + 2)
$newRunRange.LanguageID = "en-US"

# ------------------------------------------------------------------
# 5) "Gives all packets where IP is not 192.168.0.0" paragraph gets a
#    bigger left indent (720 -> 3600 twips == 36pt -> 180pt).
# ------------------------------------------------------------------
$p4 = Find-ParagraphByText $d "*Gives all packets where IP is not 192.168.0.0*"
$p4.LeftIndent = 180

Write-Output "done"
